$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 104, shifting existing rows 104:146 down to 105:147
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new record
$ws.Range("A104").Value = 3
$ws.Range("B104").Value = "Femacal de La Calera"
$ws.Range("C104").Value = "Coquimbo"
$ws.Range("D104").Value = 44466
$ws.Range("E104").Value = 5
$ws.Range("F104").Value = 100112001
$ws.Range("G104").Value = "Berenjena"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 110
$ws.Range("K104").Value = 9500
$ws.Range("L104").Value = 10000
$ws.Range("M104").Value = 9773
$ws.Range("N104").Value = "$/caja 60 unidades"
$ws.Range("O104").Value = "Región de Arica y Parinacota"
$ws.Range("P104").Value = 163
$ws.Range("Q104").Value = 60
$ws.Range("R104").Value = "Hortaliza"
